# Commit: "Added if statements and while loops to Week 2 PPT skeleton"
#
# Insert two new "skeleton" slides (title placeholder only, content
# placeholder left empty) right before the existing "What Are Python
# Packages?" slide:
#   - "The "if" Statement"
#   - "The "for" Loop"
#
# They use the same "Title and Content" layout (slideLayout2.xml /
# PowerPoint's ppLayoutText = 2) as the other skeleton slides already in
# the deck (What Are Python Packages? / Useful Package #1: NumPy / Useful
# Package #2: Pandas), and end up placed immediately ahead of that trio.

$p = $ppt.ActivePresentation

# Find where the existing "What Are Python Packages?" slide currently
# lives so the two new slides are inserted right in front of it.
$insertAt = $p.Slides.Count + 1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $title = $p.Slides.Item($i).Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "What Are Python Packages?") {
        $insertAt = $i
        break
    }
}

$ifSlide = $p.Slides.Add($insertAt, 2)
$ifSlide.Shapes.Item(1).TextFrame.TextRange.Text = "The “if” Statement"

$forSlide = $p.Slides.Add($insertAt + 1, 2)
$forSlide.Shapes.Item(1).TextFrame.TextRange.Text = "The “for” Loop"
